# Insert one new weekly record before the current row 74 ("Macroferia
# Regional de Talca - Repollo"), pushing all subsequent rows down by one
# (old row 74 -> 75, ..., old row 173 -> 174), and fill the freshly
# inserted row 74 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 74..173 down to 75..174, leaving a blank row 74 behind.
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with this week's record.
$ws.Cells.Item(74, 1).Value  = 5
$ws.Cells.Item(74, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(74, 3).Value  = 'Maule'
$ws.Cells.Item(74, 4).Value  = 44467
$ws.Cells.Item(74, 5).Value  = 7
$ws.Cells.Item(74, 6).Value  = 100112006
$ws.Cells.Item(74, 7).Value  = 'Repollo'
$ws.Cells.Item(74, 8).Value  = 'Crespo record'
$ws.Cells.Item(74, 9).Value  = 'Primera'
$ws.Cells.Item(74, 10).Value = 3000
$ws.Cells.Item(74, 11).Value = 500
$ws.Cells.Item(74, 12).Value = 500
$ws.Cells.Item(74, 13).Value = 500
$ws.Cells.Item(74, 14).Value = '$/unidad'
$ws.Cells.Item(74, 15).Value = 'Región del Maule'
$ws.Cells.Item(74, 16).Value = 500
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = 'Hortaliza'

# Keep the date column's number format consistent with the rest of the
# column (the row-insert already inherits it, but set it explicitly to
# be safe).
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
